$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data cells (rows 2-5) with new dataset values
$ws.Cells.Item(2, 1).Value = 45065.50694444445
$ws.Cells.Item(2, 2).Value = 22.58
$ws.Cells.Item(2, 3).Value = 15.542
$ws.Cells.Item(2, 4).Value = 4.221
$ws.Cells.Item(2, 5).Value = 47.493
$ws.Cells.Item(2, 6).Value = 39.284
$ws.Cells.Item(2, 7).Value = 17.769
$ws.Cells.Item(2, 8).Value = 58.8
$ws.Cells.Item(2, 9).Value = 27.341
$ws.Cells.Item(2, 10).Value = 11.61
$ws.Cells.Item(2, 11).Value = 17.881
$ws.Cells.Item(2, 12).Value = 18.828
$ws.Cells.Item(2, 13).Value = 19.728
$ws.Cells.Item(2, 14).Value = 5.673
$ws.Cells.Item(2, 15).Value = 17.67
$ws.Cells.Item(2, 16).Value = 24.849
$ws.Cells.Item(2, 17).Value = 14.79
$ws.Cells.Item(2, 18).Value = 3.779
$ws.Cells.Item(2, 19).Value = 2.46
$ws.Cells.Item(2, 20).Value = 261.617
$ws.Cells.Item(2, 21).Value = 49.202
$ws.Cells.Item(2, 22).Value = 16.31
$ws.Cells.Item(2, 23).Value = 32.642
$ws.Cells.Item(2, 24).Value = 17.025
$ws.Cells.Item(2, 25).Value = 2.109
$ws.Cells.Item(2, 26).Value = 29.266
$ws.Cells.Item(2, 27).Value = 14.407
$ws.Cells.Item(2, 28).Value = 12.944
$ws.Cells.Item(2, 29).Value = 15.145
$ws.Cells.Item(2, 30).Value = 19.485
$ws.Cells.Item(2, 31).Value = 3.64
$ws.Cells.Item(2, 32).Value = 51.902
$ws.Cells.Item(2, 33).Value = 9.071
$ws.Cells.Item(2, 34).Value = 20.391
$ws.Cells.Item(3, 1).Value = 45065.51388888889
$ws.Cells.Item(3, 2).Value = 11.05
$ws.Cells.Item(3, 3).Value = 7.614
$ws.Cells.Item(3, 4).Value = 1.637
$ws.Cells.Item(3, 5).Value = 23.392
$ws.Cells.Item(3, 6).Value = 19.364
$ws.Cells.Item(3, 7).Value = 8.696
$ws.Cells.Item(3, 8).Value = 36.481
$ws.Cells.Item(3, 9).Value = 13.38
$ws.Cells.Item(3, 10).Value = 5.702
$ws.Cells.Item(3, 11).Value = 8.653
$ws.Cells.Item(3, 12).Value = 9.42
$ws.Cells.Item(3, 13).Value = 9.766
$ws.Cells.Item(3, 14).Value = 2.78
$ws.Cells.Item(3, 15).Value = 8.647
$ws.Cells.Item(3, 16).Value = 12.154
$ws.Cells.Item(3, 17).Value = 7.483
$ws.Cells.Item(3, 18).Value = 1.578
$ws.Cells.Item(3, 19).Value = 0.902
$ws.Cells.Item(3, 20).Value = 124.311
$ws.Cells.Item(3, 21).Value = 24.305
$ws.Cells.Item(3, 22).Value = 7.982
$ws.Cells.Item(3, 23).Value = 16.011
$ws.Cells.Item(3, 24).Value = 8.590999999999999
$ws.Cells.Item(3, 25).Value = 0.976
$ws.Cells.Item(3, 26).Value = 17.178
$ws.Cells.Item(3, 27).Value = 7.05
$ws.Cells.Item(3, 28).Value = 6.448
$ws.Cells.Item(3, 29).Value = 7.535
$ws.Cells.Item(3, 30).Value = 9.771000000000001
$ws.Cells.Item(3, 31).Value = 1.294
$ws.Cells.Item(3, 32).Value = 33.032
$ws.Cells.Item(3, 33).Value = 4.38
$ws.Cells.Item(3, 34).Value = 9.978999999999999
$ws.Cells.Item(4, 1).Value = 45065.52083333334
$ws.Cells.Item(4, 2).Value = 15.374
$ws.Cells.Item(4, 3).Value = 11.097
$ws.Cells.Item(4, 4).Value = 1.289
$ws.Cells.Item(4, 5).Value = 33.024
$ws.Cells.Item(4, 6).Value = 27.321
$ws.Cells.Item(4, 7).Value = 12.098
$ws.Cells.Item(4, 8).Value = 46.426
$ws.Cells.Item(4, 9).Value = 18.615
$ws.Cells.Item(4, 10).Value = 8.146000000000001
$ws.Cells.Item(4, 11).Value = 12.239
$ws.Cells.Item(4, 12).Value = 13.335
$ws.Cells.Item(4, 13).Value = 13.92
$ws.Cells.Item(4, 14).Value = 3.864
$ws.Cells.Item(4, 15).Value = 12.031
$ws.Cells.Item(4, 16).Value = 17.03
$ws.Cells.Item(4, 17).Value = 10.243
$ws.Cells.Item(4, 18).Value = 1.074
$ws.Cells.Item(4, 19).Value = 0.78
$ws.Cells.Item(4, 20).Value = 175.792
$ws.Cells.Item(4, 21).Value = 33.62
$ws.Cells.Item(4, 22).Value = 11.105
$ws.Cells.Item(4, 23).Value = 22.437
$ws.Cells.Item(4, 24).Value = 12.007
$ws.Cells.Item(4, 25).Value = 1.45
$ws.Cells.Item(4, 26).Value = 22.398
$ws.Cells.Item(4, 27).Value = 9.808999999999999
$ws.Cells.Item(4, 28).Value = 8.795999999999999
$ws.Cells.Item(4, 29).Value = 10.317
$ws.Cells.Item(4, 30).Value = 13.926
$ws.Cells.Item(4, 31).Value = 0.784
$ws.Cells.Item(4, 32).Value = 41.814
$ws.Cells.Item(4, 33).Value = 6.197
$ws.Cells.Item(4, 34).Value = 13.884
$ws.Cells.Item(5, 1).Value = 45065.52777777778
$ws.Cells.Item(5, 2).Value = 24.02
$ws.Cells.Item(5, 3).Value = 17.71
$ws.Cells.Item(5, 4).Value = 1.36
$ws.Cells.Item(5, 5).Value = 51.94
$ws.Cells.Item(5, 6).Value = 42.96
$ws.Cells.Item(5, 7).Value = 18.9
$ws.Cells.Item(5, 8).Value = 72.16
$ws.Cells.Item(5, 9).Value = 29.09
$ws.Cells.Item(5, 10).Value = 12.91
$ws.Cells.Item(5, 11).Value = 19.31
$ws.Cells.Item(5, 12).Value = 20.94
$ws.Cells.Item(5, 13).Value = 21.98
$ws.Cells.Item(5, 14).Value = 6.04
$ws.Cells.Item(5, 15).Value = 18.8
$ws.Cells.Item(5, 16).Value = 26.76
$ws.Cells.Item(5, 17).Value = 15.81
$ws.Cells.Item(5, 18).Value = 0.9
$ws.Cells.Item(5, 19).Value = 0.93
$ws.Cells.Item(5, 20).Value = 278.82
$ws.Cells.Item(5, 21).Value = 52.55
$ws.Cells.Item(5, 22).Value = 17.35
$ws.Cells.Item(5, 23).Value = 35.34
$ws.Cells.Item(5, 24).Value = 18.75
$ws.Cells.Item(5, 25).Value = 2.36
$ws.Cells.Item(5, 26).Value = 35.22
$ws.Cells.Item(5, 27).Value = 15.33
$ws.Cells.Item(5, 28).Value = 13.59
$ws.Cells.Item(5, 29).Value = 15.97
$ws.Cells.Item(5, 30).Value = 21.94
$ws.Cells.Item(5, 31).Value = 0.5600000000000001
$ws.Cells.Item(5, 32).Value = 65.33
$ws.Cells.Item(5, 33).Value = 9.779999999999999
$ws.Cells.Item(5, 34).Value = 21.69

# Delete row 6 (dataset now has 4 rows instead of 5)
$ws.Rows.Item(6).Delete()

# Update column widths (ColumnWidth input = target stored width - 5/6)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
